$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F9").Value = 757
$ws.Range("F10").Value = 2791
$ws.Range("F11").Value = 2791
$ws.Range("F13").Value = 1845
$ws.Range("F18").Value = 6346
$ws.Range("F19").Value = 249
$ws.Range("F27").Value = 2488
$ws.Range("F29").Value = 388
$ws.Range("F35").Value = 14
$ws.Range("F36").Value = 90
$ws.Range("F39").Value = 1541
$ws.Range("F40").Value = 39
$ws.Range("F41").Value = 1493

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F11").Value = 191
$ws.Range("F18").Value = 344
$ws.Range("F19").Value = 275
$ws.Range("F20").Value = 527

# Sheet: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F7").Value = 78
$ws.Range("F8").Value = 17
$ws.Range("F9").Value = 2

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F11").Value = 78
$ws.Range("F12").Value = 78
$ws.Range("F19").Value = 2791
$ws.Range("F20").Value = 17
$ws.Range("F23").Value = 191
$ws.Range("F27").Value = 6346
$ws.Range("F28").Value = 249
$ws.Range("F34").Value = 2488
$ws.Range("F35").Value = 388
$ws.Range("F40").Value = 344
$ws.Range("F41").Value = 275
$ws.Range("F42").Value = 527
$ws.Range("F44").Value = 90
$ws.Range("F48").Value = 1541
$ws.Range("F49").Value = 39
